$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "Capstone 2"

# Clear existing contents
$ws.Cells.Clear()

# Set header row values
$ws.Range("A1").Value = "No."
$ws.Range("B1").Value = "Student Code"
$ws.Range("C1").Value = "First name"
$ws.Range("D1").Value = "Last name"
$ws.Range("E1").Value = "Group"
$ws.Range("F1").Value = "Topic"
$ws.Range("G1").Value = "Description"
